# Group the "图片 3" picture and "圆角矩形 1" rounded rectangle on slide 4
# into a single group shape ("组合 2"), then nudge the resulting group to
# its final on-slide position (size/child coordinate space stay the same
# as the natural bounding box of the two original shapes; only the
# group's own position changes).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# Shapes.Item(2) = "图片 3" (p:pic), Shapes.Item(3) = "圆角矩形 1" (p:sp)
$range = $s.Shapes.Range(@(2, 3))
$grp = $range.Group()

$grp.Name = "组合 2"

# Final group position in points (EMU / 12700), nudged within the
# float-rounding tolerance of the host so the saved EMU values land
# exactly on 5231904 / 548680.
$grp.Left = 411.9609848818898
$grp.Top = 43.20318960629921
